# Auto-generated edit script: apply cached numeric updates scraped by the scheduled runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 405.75
$ws.Range("I28").Value = 509.18182
$ws.Range("K28").Value = 509.18182
$ws.Range("M28").Value = -24.18182000000002
$ws.Range("H48").Value = 1871
$ws.Range("J48").Value = 1833.3334
$ws.Range("L48").Value = 5500.0002
$ws.Range("N48").Value = -6084.0002
$ws.Range("H56").Value = 1871
$ws.Range("J56").Value = 1833.3334
$ws.Range("L56").Value = 5500.0002
$ws.Range("N56").Value = -6568.0002
$ws.Range("H70").Value = 2293.0833
$ws.Range("I70").Value = 2478
$ws.Range("J70").Value = 2161
$ws.Range("K70").Value = 7434
$ws.Range("L70").Value = 6483
$ws.Range("M70").Value = -7164
$ws.Range("N70").Value = -7023
$ws.Range("H73").Value = 2293.0833
$ws.Range("I73").Value = 2478
$ws.Range("J73").Value = 2161
$ws.Range("K73").Value = 7434
$ws.Range("L73").Value = 6483
$ws.Range("M73").Value = -6498
$ws.Range("N73").Value = -8355
$ws.Range("H88").Value = 4367.5
$ws.Range("I88").Value = 4656.6665
$ws.Range("K88").Value = 4656.6665
$ws.Range("M88").Value = -4250.6665
$ws.Range("H91").Value = 4367.5
$ws.Range("I91").Value = 4656.6665
$ws.Range("K91").Value = 4656.6665
$ws.Range("M91").Value = -3252.6665
$ws.Range("H140").Value = 99874.5
$ws.Range("J140").Value = 99874.5
$ws.Range("L140").Value = 99874.5
$ws.Range("N140").Value = -110234.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 658.6667
$ws.Range("I4").Value = 390
$ws.Range("J4").Value = 2002
$ws.Range("K4").Value = 390
$ws.Range("L4").Value = 2002
$ws.Range("M4").Value = -274
$ws.Range("N4").Value = -2234
$ws.Range("H5").Value = 2304.875
$ws.Range("I5").Value = 110
$ws.Range("J5").Value = 4499.75
$ws.Range("K5").Value = 110
$ws.Range("L5").Value = 4499.75
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = -4723.75
$ws.Range("H88").Value = 2484.652
$ws.Range("I88").Value = 2077.7144
$ws.Range("J88").Value = 2662.6875
$ws.Range("K88").Value = 2077.7144
$ws.Range("L88").Value = 2662.6875
$ws.Range("M88").Value = -1671.7144
$ws.Range("N88").Value = -3474.6875
$ws.Range("H91").Value = 2484.652
$ws.Range("I91").Value = 2077.7144
$ws.Range("J91").Value = 2662.6875
$ws.Range("K91").Value = 2077.7144
$ws.Range("L91").Value = 2662.6875
$ws.Range("M91").Value = -673.7143999999998
$ws.Range("N91").Value = -5470.6875
$ws.Range("H132").Value = 738714.75
$ws.Range("I132").Value = 809977.6
$ws.Range("K132").Value = 2429932.8
$ws.Range("M132").Value = -2427402.8
$ws.Range("H134").Value = 54982.332
$ws.Range("J134").Value = 54982.332
$ws.Range("L134").Value = 54982.332
$ws.Range("N134").Value = -65122.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2304.875
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 4499.75
$ws.Range("K4").Value = 110
$ws.Range("L4").Value = 4499.75
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = -4729.75
$ws.Range("H20").Value = 47782.39
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2494
$ws.Range("H86").Value = 538.1
$ws.Range("I86").Value = 486.77777
$ws.Range("K86").Value = 486.77777
$ws.Range("M86").Value = 636.2222300000001
$ws.Range("H89").Value = 538.1
$ws.Range("I89").Value = 486.77777
$ws.Range("K89").Value = 2433.88885
$ws.Range("M89").Value = 3182.11115
$ws.Range("H94").Value = 1437.5454
$ws.Range("I94").Value = 1449.2307
$ws.Range("J94").Value = 1394.1428
$ws.Range("K94").Value = 1449.2307
$ws.Range("L94").Value = 1394.1428
$ws.Range("M94").Value = -998.2307000000001
$ws.Range("N94").Value = -2296.1428
$ws.Range("H105").Value = 11099.8
$ws.Range("I105").Value = 12500
$ws.Range("J105").Value = 8999.5
$ws.Range("K105").Value = 12500
$ws.Range("L105").Value = 8999.5
$ws.Range("M105").Value = -10753
$ws.Range("N105").Value = -12493.5
$ws.Range("H107").Value = 1550
$ws.Range("I107").Value = 1066.6666
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1066.6666
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 853.3334
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 147
$ws.Range("I2").Value = 147
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 147
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H82").Value = 37900
$ws.Range("J82").Value = 37900
$ws.Range("L82").Value = 37900
$ws.Range("N82").Value = -38622
$ws.Range("H85").Value = 37900
$ws.Range("J85").Value = 37900
$ws.Range("L85").Value = 37900
$ws.Range("N85").Value = -40396
$ws.Range("H141").Value = 187422.42
$ws.Range("J141").Value = 210327.08
$ws.Range("L141").Value = 210327.08
$ws.Range("N141").Value = -220687.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2034621.2
$ws.Range("I5").Value = 1984529.2
$ws.Range("K5").Value = 5953587.6
$ws.Range("M5").Value = -5953475.6
$ws.Range("H46").Value = 2188.56
$ws.Range("I46").Value = 1817.9445
$ws.Range("K46").Value = 5453.833500000001
$ws.Range("M46").Value = -5362.833500000001
$ws.Range("H107").Value = 3390.0625
$ws.Range("I107").Value = 444.6
$ws.Range("J107").Value = 4728.909
$ws.Range("K107").Value = 1333.8
$ws.Range("L107").Value = 14186.727
$ws.Range("M107").Value = 586.1999999999998
$ws.Range("N107").Value = -18026.727
$ws.Range("H131").Value = 4237.143
$ws.Range("J131").Value = 5387.1333
$ws.Range("L131").Value = 16161.3999
$ws.Range("N131").Value = -26241.3999
$ws.Range("H135").Value = 2034621.2
$ws.Range("I135").Value = 1984529.2
$ws.Range("K135").Value = 17860762.8
$ws.Range("M135").Value = -17858227.8
$ws.Range("H139").Value = 8284.9
$ws.Range("I139").Value = 5658.778
$ws.Range("J139").Value = 10433.546
$ws.Range("K139").Value = 16976.334
$ws.Range("L139").Value = 31300.638
$ws.Range("M139").Value = -11836.334
$ws.Range("N139").Value = -41580.638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31313.781
$ws.Range("I2").Value = 38503.117
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 38503.117
$ws.Range("L2").Value = 160
$ws.Range("M2").Value = -38390.117
$ws.Range("N2").Value = -386
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = ""
$ws.Range("H80").Value = 8026
$ws.Range("I80").Value = 5044.2856
$ws.Range("J80").Value = 14983.333
$ws.Range("K80").Value = 5044.2856
$ws.Range("L80").Value = 14983.333
$ws.Range("M80").Value = -4046.2856
$ws.Range("N80").Value = -16979.333
$ws.Range("H83").Value = 8026
$ws.Range("I83").Value = 5044.2856
$ws.Range("J83").Value = 14983.333
$ws.Range("K83").Value = 25221.428
$ws.Range("L83").Value = 74916.66500000001
$ws.Range("M83").Value = -20229.428
$ws.Range("N83").Value = -84900.66500000001
$ws.Range("H97").Value = 1405.0238
$ws.Range("I97").Value = 1374.4546
$ws.Range("K97").Value = 1374.4546
$ws.Range("M97").Value = -878.4546
$ws.Range("H113").Value = 3991.111
$ws.Range("I113").Value = 3485.6
$ws.Range("J113").Value = 4623
$ws.Range("K113").Value = 3485.6
$ws.Range("L113").Value = 4623
$ws.Range("M113").Value = -1315.6
$ws.Range("N113").Value = -8963
$ws.Range("H132").Value = 50586.832
$ws.Range("J132").Value = 41540.332
$ws.Range("L132").Value = 124620.996
$ws.Range("N132").Value = -129680.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12469.647
$ws.Range("I7").Value = 29465.834
$ws.Range("K7").Value = 29465.834
$ws.Range("M7").Value = -29353.834
$ws.Range("H25").Value = 8998.666999999999
$ws.Range("I25").Value = 8998
$ws.Range("K25").Value = 8998
$ws.Range("M25").Value = -8768
$ws.Range("H126").Value = 12469.647
$ws.Range("I126").Value = 29465.834
$ws.Range("K126").Value = 88397.50199999999
$ws.Range("M126").Value = -85927.50199999999
$ws.Range("H132").Value = 12986374
$ws.Range("I132").Value = 29214598
$ws.Range("J132").Value = 3795
$ws.Range("K132").Value = 87643794
$ws.Range("L132").Value = 11385
$ws.Range("M132").Value = -87641264
$ws.Range("N132").Value = -16445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4039.5
$ws.Range("I81").Value = 1460
$ws.Range("J81").Value = 5882
$ws.Range("K81").Value = 2920
$ws.Range("L81").Value = 11764
$ws.Range("M81").Value = -1859
$ws.Range("N81").Value = -13886
$ws.Range("H84").Value = 4039.5
$ws.Range("I84").Value = 1460
$ws.Range("J84").Value = 5882
$ws.Range("K84").Value = 14600
$ws.Range("L84").Value = 58820
$ws.Range("M84").Value = -9296
$ws.Range("N84").Value = -69428
$ws.Range("H132").Value = 6174985.5
$ws.Range("J132").Value = 3366.6667
$ws.Range("L132").Value = 10100.0001
$ws.Range("N132").Value = -15160.0001
$ws.Range("H136").Value = 15219476
$ws.Range("I136").Value = 3624726
$ws.Range("K136").Value = 10874178
$ws.Range("M136").Value = -10871628
